$d = $word.ActiveDocument

# WdColorIndex.wdBrightGreen = 4  -> OOXML <w:highlight w:val="green"/>
$wdBrightGreen = 4

# Paragraphs (1-based, matching $d.Paragraphs) that become highlighted green
# in this revision (the ones marked "confirmed" complete by the author):
#  10 - "(2.5 points) ... seeded data ..."
#  11 - "(5 points) ... Postman ..."                (already green; left as-is)
#  12 - "(10 points) ... aesthetically pleasing ..."
#  13 - "(5 points) ... see the details of a movie ..."
#  15 - "(5 points) ... add a new movie ..."
# Paragraph 14 ("... update the details ...") stays unhighlighted.

$targets = @(10, 11, 12, 13, 15)

foreach ($idx in $targets) {
    $p = $d.Paragraphs($idx)
    $p.Range.HighlightColorIndex = $wdBrightGreen
}
